# Updated cryptos list values (Price + Volume(1h)) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Force the cell to be stored as text even when the string looks numeric
    # (e.g. "236.11"), then restore the default "Normal" style so no stray
    # number-format style is left behind on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.222.72"
$ws.Range("E2").Value = "  +0.32%  "
Set-TextValue $ws.Range("D3") "1.859.72"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  -0.05%  "
Set-TextValue $ws.Range("D5") "236.11"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("E7").Value = "  +0.40%  "
Set-TextValue $ws.Range("D8") "0.2870"
$ws.Range("E8").Value = "  +1.77%  "
Set-TextValue $ws.Range("D9") "0.06519"
$ws.Range("E9").Value = "  -0.59%  "
Set-TextValue $ws.Range("D10") "21.66"
$ws.Range("E10").Value = "  +8.21%  "
Set-TextValue $ws.Range("D11") "0.07895"
$ws.Range("E11").Value = "  +0.91%  "
Set-TextValue $ws.Range("D12") "97.34"
$ws.Range("E12").Value = "  +0.56%  "
Set-TextValue $ws.Range("D13") "1.864.96"
$ws.Range("E13").Value = "  +0.37%  "
Set-TextValue $ws.Range("D14") "5.157"
$ws.Range("E14").Value = "  +0.92%  "
Set-TextValue $ws.Range("D15") "0.6790"
$ws.Range("E15").Value = "  +1.99%  "
Set-TextValue $ws.Range("D16") "279.24"
$ws.Range("E16").Value = "  -1.09%  "
Set-TextValue $ws.Range("D17") "30.219.18"
$ws.Range("E17").Value = "  +0.18%  "
Set-TextValue $ws.Range("D18") "13.48"
$ws.Range("E18").Value = "  +6.90%  "
Set-TextValue $ws.Range("D19") "1.0000"
$ws.Range("E19").Value = "  -0.04%  "
Set-TextValue $ws.Range("D20") "5.372"
Set-TextValue $ws.Range("D21") "2.110.00"
$ws.Range("E21").Value = "  +0.12%  "
Set-TextValue $ws.Range("D22") "0.000007297"
$ws.Range("E22").Value = "  +0.94%  "
$ws.Range("E23").Value = "  -0.02%  "
Set-TextValue $ws.Range("D24") "6.161"
$ws.Range("E24").Value = "  +0.40%  "
Set-TextValue $ws.Range("D25") "167.19"
$ws.Range("E25").Value = "  -0.42%  "
Set-TextValue $ws.Range("D26") "9.216"
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("E27").Value = "  +1.01%  "
Set-TextValue $ws.Range("D28") "1.930"
$ws.Range("E28").Value = "  +0.84%  "
Set-TextValue $ws.Range("D29") "1.381"
$ws.Range("E29").Value = "  +3.54%  "
Set-TextValue $ws.Range("D30") "0.09718"
$ws.Range("E30").Value = "  +1.59%  "
Set-TextValue $ws.Range("D31") "4.365"
$ws.Range("E31").Value = "  -1.18%  "
Set-TextValue $ws.Range("D32") "1.479"
$ws.Range("E32").Value = "  +0.54%  "
Set-TextValue $ws.Range("D33") "4.043"
$ws.Range("E33").Value = "  -1.44%  "
Set-TextValue $ws.Range("D34") "0.04722"
$ws.Range("E34").Value = "  +1.38%  "
Set-TextValue $ws.Range("D35") "1.131"
$ws.Range("E35").Value = "  +2.96%  "
Set-TextValue $ws.Range("D36") "0.7068"
$ws.Range("E36").Value = "  +0.76%  "
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("E38").Value = "  +0.73%  "
Set-TextValue $ws.Range("D39") "2.631"
$ws.Range("E39").Value = "  +4.81%  "
Set-TextValue $ws.Range("D40") "6.331"
$ws.Range("E40").Value = "  +0.15%  "
Set-TextValue $ws.Range("D41") "74.46"
$ws.Range("E41").Value = "  +3.16%  "
Set-TextValue $ws.Range("D42") "1.949"
$ws.Range("E42").Value = "  +1.18%  "
Set-TextValue $ws.Range("D43") "0.8476"
$ws.Range("E43").Value = "  -0.50%  "
Set-TextValue $ws.Range("D44") "0.4171"
$ws.Range("E44").Value = "  +0.52%  "
Set-TextValue $ws.Range("D45") "0.9996"
Set-TextValue $ws.Range("D46") "103.31"
$ws.Range("E46").Value = "  -0.55%  "
Set-TextValue $ws.Range("D47") "968.68"
$ws.Range("E47").Value = "  -2.11%  "
Set-TextValue $ws.Range("D48") "7.176"
$ws.Range("E48").Value = "  -0.80%  "
Set-TextValue $ws.Range("D49") "9.232"
$ws.Range("E49").Value = "  +0.34%  "
$ws.Range("E50").Value = "  +0.31%  "
Set-TextValue $ws.Range("D51") "0.05638"
$ws.Range("E51").Value = "  +0.20%  "
